$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a second copy of the "2BJs/4BJs" summary block (originally at rows
# 38-45, columns B:C) further down the sheet as a combined 4-column block at
# rows 52-55, so it reads alongside its SE markers on the same rows.

# Header row: labels + "SE" markers
$ws.Range("B52").Value = "2BJs 2Machines, 8 replicas, 32 exchanges"
$ws.Range("C52").Value = "SE"
$ws.Range("D52").Value = "4BJs, 4Machines - 16 replicas/64 exchanges"
$ws.Range("E52").Value = "SE"

# Synchronous row
$ws.Range("A53").Value = "Synchronous "
$ws.Range("B53").Value = 805
$ws.Range("C53").Value = 14.5
$ws.Range("D53").Value = 1179.8
$ws.Range("E53").Value = 10.05

# Asynchronous - Centralized row
$ws.Range("A54").Value = "Asynchronous - Centralized"
$ws.Range("B54").Value = 632
$ws.Range("C54").Value = 7.11
$ws.Range("D54").Value = 685
$ws.Range("E54").Value = 5.57

# Asynchronous - Decentralized row
$ws.Range("A55").Value = "Asynchronous - Decentralized"
$ws.Range("B55").Value = 607.8
$ws.Range("C55").Value = 1.66
$ws.Range("D55").Value = 641
$ws.Range("E55").Value = 9.17

# Scroll the view down and move the selection the way the author left it.
$excel.ActiveWindow.ScrollRow = 26
$ws.Range("A56").Select()
